$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 4432.8887
$ws.Range("I32").Value = 949.5
$ws.Range("J32").Value = 5428.143
$ws.Range("K32").Value = 949.5
$ws.Range("L32").Value = 5428.143
$ws.Range("M32").Value = -623.5
$ws.Range("N32").Value = -6080.143
$ws.Range("H55").Value = 332.375
$ws.Range("I55").Value = 328.8
$ws.Range("J55").Value = 338.33334
$ws.Range("K55").Value = 328.8
$ws.Range("L55").Value = 338.33334
$ws.Range("M55").Value = -114.8
$ws.Range("N55").Value = -766.33334
$ws.Range("H106").Value = 6444.2
$ws.Range("I106").Value = 11111
$ws.Range("J106").Value = 3333
$ws.Range("K106").Value = 11111
$ws.Range("L106").Value = 3333
$ws.Range("M106").Value = -10480
$ws.Range("N106").Value = -4595
$ws.Range("H125").Value = 3038705
$ws.Range("I125").Value = 4139693.8
$ws.Range("K125").Value = 37257244.2
$ws.Range("M125").Value = -37254784.2
$ws.Range("H130").Value = 250000
$ws.Range("J130").Value = 250000
$ws.Range("L130").Value = 250000
$ws.Range("N130").Value = -260040
$ws.Range("H131").Value = 2822.4
$ws.Range("I131").Value = 2247.6667
$ws.Range("K131").Value = 6743.000100000001
$ws.Range("M131").Value = -1703.000100000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 562840.8
$ws.Range("I61").Value = 2294.7441
$ws.Range("K61").Value = 2294.7441
$ws.Range("M61").Value = -2082.7441
$ws.Range("H132").Value = 6023577.5
$ws.Range("I132").Value = 3154.077
$ws.Range("J132").Value = 14719745
$ws.Range("K132").Value = 9462.231
$ws.Range("L132").Value = 44159235
$ws.Range("M132").Value = -6932.231
$ws.Range("N132").Value = -44164295
$ws.Range("H136").Value = 562840.8
$ws.Range("I136").Value = 2294.7441
$ws.Range("K136").Value = 6884.2323
$ws.Range("M136").Value = -4334.2323

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 43470.5
$ws.Range("I99").Value = 40926.25
$ws.Range("J99").Value = 53647.5
$ws.Range("K99").Value = 40926.25
$ws.Range("L99").Value = 53647.5
$ws.Range("M99").Value = -39428.25
$ws.Range("N99").Value = -56643.5
$ws.Range("H107").Value = 950
$ws.Range("I107").Value = 950
$ws.Range("K107").Value = 950
$ws.Range("M107").Value = 970
$ws.Range("H134").Value = 9101.5
$ws.Range("I134").Value = 3293.1562
$ws.Range("J134").Value = 20718.188
$ws.Range("K134").Value = 9879.4686
$ws.Range("L134").Value = 62154.564
$ws.Range("M134").Value = -7344.4686
$ws.Range("N134").Value = -67224.564

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5972.93
$ws.Range("I31").Value = 1277.1282
$ws.Range("J31").Value = 51757
$ws.Range("K31").Value = 1277.1282
$ws.Range("L31").Value = 51757
$ws.Range("M31").Value = -982.1282000000001
$ws.Range("N31").Value = -52347
$ws.Range("H34").Value = 5972.93
$ws.Range("I34").Value = 1277.1282
$ws.Range("J34").Value = 51757
$ws.Range("K34").Value = 1277.1282
$ws.Range("L34").Value = 51757
$ws.Range("M34").Value = -1075.1282
$ws.Range("N34").Value = -52161
$ws.Range("H58").Value = 15719.518
$ws.Range("I58").Value = 6409.6665
$ws.Range("J58").Value = 25694.357
$ws.Range("K58").Value = 6409.6665
$ws.Range("L58").Value = 25694.357
$ws.Range("M58").Value = -6206.6665
$ws.Range("N58").Value = -26100.357
$ws.Range("H80").Value = 18875
$ws.Range("J80").Value = 14333.333
$ws.Range("L80").Value = 14333.333
$ws.Range("N80").Value = -16579.333
$ws.Range("H83").Value = 18875
$ws.Range("J83").Value = 14333.333
$ws.Range("L83").Value = 42999.999
$ws.Range("N83").Value = -54231.999
$ws.Range("H107").Value = 2443.2727
$ws.Range("I107").Value = 2230.4285
$ws.Range("J107").Value = 2815.75
$ws.Range("K107").Value = 2230.4285
$ws.Range("L107").Value = 2815.75
$ws.Range("M107").Value = -310.4285
$ws.Range("N107").Value = -6655.75
$ws.Range("H134").Value = 27783240
$ws.Range("I134").Value = 2231.0715
$ws.Range("J134").Value = 125016776
$ws.Range("K134").Value = 6693.2145
$ws.Range("L134").Value = 375050328
$ws.Range("M134").Value = -4158.2145
$ws.Range("N134").Value = -375055398
$ws.Range("H136").Value = 15719.518
$ws.Range("I136").Value = 6409.6665
$ws.Range("J136").Value = 25694.357
$ws.Range("K136").Value = 19228.9995
$ws.Range("L136").Value = 77083.071
$ws.Range("M136").Value = -16678.9995
$ws.Range("N136").Value = -82183.071

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 991.6667
$ws.Range("I5").Value = 1026.5294
$ws.Range("J5").Value = 907
$ws.Range("K5").Value = 3079.5882
$ws.Range("L5").Value = 2721
$ws.Range("M5").Value = -2967.5882
$ws.Range("N5").Value = -2945
$ws.Range("H38").Value = 198.88889
$ws.Range("I38").Value = 170
$ws.Range("K38").Value = 510
$ws.Range("M38").Value = -163
$ws.Range("H135").Value = 991.6667
$ws.Range("I135").Value = 1026.5294
$ws.Range("J135").Value = 907
$ws.Range("K135").Value = 9238.764599999999
$ws.Range("L135").Value = 8163
$ws.Range("M135").Value = -6703.764599999999
$ws.Range("N135").Value = -13233

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 5654144
$ws.Range("I126").Value = 3595578.8
$ws.Range("K126").Value = 10786736.4
$ws.Range("M126").Value = -10784266.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 52632580
$ws.Range("I22").Value = 904.5454999999999
$ws.Range("J22").Value = 125001140
$ws.Range("K22").Value = 904.5454999999999
$ws.Range("L22").Value = 125001140
$ws.Range("M22").Value = -609.5454999999999
$ws.Range("N22").Value = -125001730
$ws.Range("H27").Value = 52632580
$ws.Range("I27").Value = 904.5454999999999
$ws.Range("J27").Value = 125001140
$ws.Range("K27").Value = 904.5454999999999
$ws.Range("L27").Value = 125001140
$ws.Range("M27").Value = -797.5454999999999
$ws.Range("N27").Value = -125001354
$ws.Range("H55").Value = 1859.6136
$ws.Range("I55").Value = 1938.0588
$ws.Range("J55").Value = 1810.2222
$ws.Range("K55").Value = 1938.0588
$ws.Range("L55").Value = 1810.2222
$ws.Range("M55").Value = -1765.0588
$ws.Range("N55").Value = -2156.2222
$ws.Range("H130").Value = 30000
$ws.Range("J130").Value = 30000
$ws.Range("L130").Value = 30000
$ws.Range("N130").Value = -40040

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 824.3889
$ws.Range("I107").Value = 934.9286
$ws.Range("K107").Value = 2804.7858
$ws.Range("M107").Value = -884.7857999999997
$ws.Range("H126").Value = 5577456
$ws.Range("I126").Value = 24560.75
$ws.Range("K126").Value = 73682.25
$ws.Range("M126").Value = -71212.25
$ws.Range("H132").Value = 589074.0600000001
$ws.Range("I132").Value = 8083.7856
$ws.Range("J132").Value = 2215846.8
$ws.Range("K132").Value = 24251.3568
$ws.Range("L132").Value = 6647540.399999999
$ws.Range("M132").Value = -21721.3568
$ws.Range("N132").Value = -6652600.399999999
$ws.Range("H136").Value = 373189.7
$ws.Range("I136").Value = 1588.0454
$ws.Range("K136").Value = 4764.1362
$ws.Range("M136").Value = -2214.1362
